$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "484.91").
# Excel would auto-convert such text to a numeric value on assignment, which
# would not match the source data (stored as literal text, e.g. "0.130" keeps
# its trailing zero, "69.055.81" has two dots, etc). Force these specific cells
# to Text format first so the assigned string is preserved verbatim.
$textFormatRows = 4,5,6,7,9,11,12,14,16,18,21,22,23,24,25,26,27,28,29,30,31,32,33,35,36,38,40,42,49,50
foreach ($r in $textFormatRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '69.055.81'
$ws.Range('E2').Value = '  +2.65%  '

# Row 3
$ws.Range('D3').Value = '3.945.11'
$ws.Range('E3').Value = '  +0.98%  '

# Row 4
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '

# Row 5
$ws.Range('D5').Value = '484.91'
$ws.Range('E5').Value = '  +3.38%  '

# Row 6
$ws.Range('D6').Value = '146.63'
$ws.Range('E6').Value = '  +0.31%  '

# Row 7
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -1.76%  '

# Row 8
$ws.Range('E8').Value = '  -0.11%  '

# Row 9
$ws.Range('D9').Value = '0.725'
$ws.Range('E9').Value = '  -2.47%  '

# Row 10
$ws.Range('E10').Value = '  +8.27%  '

# Row 11
$ws.Range('D11').Value = '0.0000358'
$ws.Range('E11').Value = '  +13.64%  '

# Row 12
$ws.Range('D12').Value = '42.71'
$ws.Range('E12').Value = '  -1.87%  '

# Row 13
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.565.49'
$ws.Range('E13').Value = '  -0.09%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '10.48'
$ws.Range('E14').Value = '  +0.85%  '

# Row 15
$ws.Range('D15').Value = '3.973.37'
$ws.Range('E15').Value = '  +1.18%  '

# Row 16
$ws.Range('D16').Value = '14.61'
$ws.Range('E16').Value = '  -1.66%  '

# Row 17
$ws.Range('E17').Value = '  -0.24%  '

# Row 18
$ws.Range('D18').Value = '19.71'
$ws.Range('E18').Value = '  -1.74%  '

# Row 19
$ws.Range('E19').Value = '  -2.77%  '

# Row 20
$ws.Range('D20').Value = '69.079.19'
$ws.Range('E20').Value = '  +2.08%  '

# Row 21
$ws.Range('D21').Value = '435.98'
$ws.Range('E21').Value = '  +1.03%  '

# Row 22
$ws.Range('D22').Value = '14.61'
$ws.Range('E22').Value = '  -0.87%  '

# Row 23
$ws.Range('D23').Value = '3.36'
$ws.Range('E23').Value = '  +2.64%  '

# Row 24
$ws.Range('D24').Value = '87.90'
$ws.Range('E24').Value = '  -1.09%  '

# Row 25
$ws.Range('D25').Value = '11.70'
$ws.Range('E25').Value = '  +16.48%  '

# Row 26
$ws.Range('D26').Value = '3.58'
$ws.Range('E26').Value = '  -0.59%  '

# Row 27
$ws.Range('D27').Value = '38.43'
$ws.Range('E27').Value = '  +1.02%  '

# Row 28
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').Value = '5.90'
$ws.Range('E28').Value = '  +7.64%  '

# Row 29
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '10.25'
$ws.Range('E29').Value = '  +0.42%  '

# Row 30
$ws.Range('D30').Value = '712.63'
$ws.Range('E30').Value = '  -2.46%  '

# Row 31
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').Value = '13.28'
$ws.Range('E31').Value = '  -3.17%  '

# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.130'
$ws.Range('E32').Value = '  -3.62%  '

# Row 33
$ws.Range('D33').Value = '2.86'
$ws.Range('E33').Value = '  +4.50%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0937'
$ws.Range('E34').Value = '  +36.10%  '

# Row 35
$ws.Range('D35').Value = '41.34'
$ws.Range('E35').Value = '  -4.66%  '

# Row 36
$ws.Range('D36').Value = '58.76'
$ws.Range('E36').Value = '  +2.35%  '

# Row 37
$ws.Range('E37').Value = '  -7.07%  '

# Row 38
$ws.Range('D38').Value = '5.64'
$ws.Range('E38').Value = '  +2.72%  '

# Row 39
$ws.Range('E39').Value = '  -0.11%  '

# Row 40
$ws.Range('D40').Value = '0.0472'
$ws.Range('E40').Value = '  -1.87%  '

# Row 41
$ws.Range('E41').Value = '  +8.65%  '

# Row 42
$ws.Range('D42').Value = '3.05'
$ws.Range('E42').Value = '  +9.98%  '

# Row 43
$ws.Range('E43').Value = '  +2.80%  '

# Row 44
$ws.Range('E44').Value = '  -1.40%  '

# Row 45
$ws.Range('E45').Value = '  +0.09%  '

# Row 46
$ws.Range('E46').Value = '  -0.20%  '

# Row 47
$ws.Range('E47').Value = '  -2.17%  '

# Row 48
$ws.Range('E48').Value = '  -1.09%  '

# Row 49
$ws.Range('D49').Value = '148.16'
$ws.Range('E49').Value = '  +2.76%  '

# Row 50
$ws.Range('D50').Value = '3.11'
$ws.Range('E50').Value = '  -4.57%  '

# Row 51
$ws.Range('E51').Value = '  -1.57%  '
